# Generate Report for Handoff
# Updates the localization-status workbook after a fresh handoff of the
# f785546e-6d87-42af-8ccc-5203ee0dac32.md file (row 7 in each sheet).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: refresh "Latest HO Xliff Generate Date" for the f785546e row
$wsOverview.Range("G7").Value = "2016-09-01 18:48:55"

# zh-cn sheet: refresh handoff info for the f785546e row
$wsZhCn.Range("G7").Value = "2016-09-01 18:48:50"
$wsZhCn.Range("H7").Value = "f785546e-6d87-42af-8ccc-5203ee0dac32.cdfaf3490c42c8626702b019a8d20abff4381555.zh-cn.xlf"

# de-de sheet: refresh "Latest Handoff Datetime" for the f785546e row
$wsDeDe.Range("H7").Value = "2016-09-01 18:48:55"
